# "deep learning notes.docx" update:
#   - after the existing "Optimizer in neural network" heading, add an
#     explanatory paragraph, a "Types of Optimizer" Heading2 with a bulleted
#     (numId 13) list of eight optimizers, and a new
#     "Customer Churn Prediction using ANN" Heading1.
#   - the hidden "_GoBack" bookmark, which currently sits inside the
#     "Optimizer in neural network" paragraph, ends up in its own empty
#     paragraph at the very end of the document body (after all of the new
#     content), exactly as it was left by the editing session that produced
#     the target revision.

$d = $word.ActiveDocument

# Locate the paragraph that currently holds both the
# "Optimizer in neural network" heading text and the "_GoBack" bookmark -
# it is the last paragraph in the document body.
$count = $d.Paragraphs.Count
$targetPara = $d.Paragraphs.Item($count)

if ($targetPara.Range.Text -notmatch "Optimizer in neural network") {
    throw "Could not locate the 'Optimizer in neural network' paragraph"
}

# Insert the new content right after the heading's text, i.e. immediately
# before the bookmark markup that currently closes out the paragraph. That
# keeps the bookmark untouched for now and lands every paragraph of the new
# fragment - including the bookmark-holder paragraph at the very end - as a
# clean, freestanding paragraph.
$insertPos = $targetPara.Range.End - 1
$insertionPoint = $d.Range($insertPos, $insertPos)

$fragment = '<w:p><w:r><w:t>Optimizers are algorithms or methods used to change the attributes of your neural network such as weights and learning rate in order to reduce the losses.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Types of Optimizer</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr></w:pPr><w:r><w:t>Gradient Descent</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr></w:pPr><w:r><w:t>Stochastic Gradient Descent</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr></w:pPr><w:r><w:t>Stochastic Gradient Descent with momentum</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr></w:pPr><w:r><w:t>Mini-Batch Gradient Descent</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr></w:pPr><w:r><w:t>Adagrad</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr></w:pPr><w:r><w:t>RMSProp</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr></w:pPr><w:r><w:t>AdaDelta</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="13"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Adam </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/><w:t>(mostly use this)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Customer Churn Prediction using ANN</w:t></w:r></w:p><w:p><w:bookmarkStart w:id="99" w:name="_GoBack"/><w:bookmarkEnd w:id="99"/></w:p>'

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $fragment + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($xml) | Out-Null

# Remove the original "_GoBack" bookmark - the one now wrapping the new
# trailing paragraph (added above) is the one that survives.
$d.Bookmarks.Item("_GoBack").Delete()
